$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 38: columns G:J get the same labels as C22:F22 (shared strings 0..3)
$ws.Range("G38").Value = "HK-2"
$ws.Range("H38").Value = "UMRC6"
$ws.Range("I38").Value = "UOK262"
$ws.Range("J38").Value = "UOK + DIDS"

# Row 39: label + averages
$ws.Range("F39").Value = "Flow_Lac"
$ws.Range("G39").Formula = "=AVERAGE(F`$1:F`$3)"
$ws.Range("H39").Formula = "=AVERAGE(F`$4:F`$6)"
$ws.Range("I39").Formula = "=AVERAGE(F`$9:F`$11)"
$ws.Range("J39").Formula = "=AVERAGE(F`$13:F`$16)"

# Row 40: standard error
$ws.Range("G40").Formula = "=STDEV(F`$1:F`$3)/SQRT(COUNT(F`$1:F`$3))"
$ws.Range("H40").Formula = "=STDEV(F`$4:F`$6)/SQRT(COUNT(F`$4:F`$6))"
$ws.Range("I40").Formula = "=STDEV(F`$9:F`$11)/SQRT(COUNT(F`$9:F`$11))"
$ws.Range("J40").Formula = "=STDEV(F`$13:F`$16)/SQRT(COUNT(F`$13:F`$16))"

# Update view selection/scroll like diff indicates
$ws.Range("F38:J40").Select()
